# cv122053a.xlsx - "correção nos dados e inicio da analise PNAD 2009"
#
# The sheet had two "label only" rows (row 5 "situação do domicílio" and
# row 8 "grandes regiões e unidades da federação") that carried no data
# values of their own - they were just section headers squeezed into the
# data grid. This edit removes those two rows (so the rows below shift up
# and keep their own values), and fixes the second header row so that the
# "total" column header is repeated in B2/C2/F2 instead of the stray
# "unnamed: 1_level_1" / "unnamed: 5_level_1" placeholder labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix header row 2 -----------------------------------------------------
# C2 already holds the correct "total" label; reuse its text for B2 and F2,
# which currently hold the placeholder "unnamed: 1_level_1" / "unnamed:
# 5_level_1" strings.
$totalText = $ws.Range("C2").Value2
$ws.Range("B2").Value = $totalText
$ws.Range("F2").Value = $totalText

# --- Remove the two label-only rows ---------------------------------------
# Row 5 ("situação do domicílio") has no data; deleting it shifts every
# subsequent row up by one.
$ws.Rows(5).Delete()

# After the row-5 deletion, the old row 8 ("grandes regiões e unidades da
# federação") is now row 7; it also has no data, so remove it too.
$ws.Rows(7).Delete()
